$wb = $excel.ActiveWorkbook

# Rename the data sheet to the template name (drop the "naming convention" prefix).
$ws1 = $wb.Worksheets.Item("3ASY01_RNASeq")
$ws2 = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws1.Name = "rnaseq_assay"

# Update the metadata "Name" field to match the new template name.
$ws2.Range("B2").Value = "RNASeq Assay"

# Update selections: metadata sheet selection moves to B3 (no longer active tab) ...
$ws2.Range("B3").Select()

# ... and the renamed data sheet becomes the active tab with selection at B40.
$ws1.Activate()
$ws1.Range("B40").Select()
